# Apply edits to "DummyJson API Testing.xlsx" as described by the commit
# "Update Readme.MD and POJO classes":
#   - Rename sheet "Cart and Order Data Retrieval" -> "Cart Data Retrieval"
#   - Split the single shared string "Fetch Single User" into two distinct
#     test-case titles in the "User Data Retrieval" sheet:
#       B3: "Fetch Single User" -> "Fetch Valid Single User"
#       B4: "Fetch Single User" -> "Fetch Invalid Single User"
#   - Update the saved selections/active sheet to reflect where the author
#     was last working (User Data Retrieval sheet active with J9 selected,
#     Cart Data Retrieval selection at D22, Recipe Data Retrieval selection
#     at D31).

$wb = $excel.ActiveWorkbook

# --- Rename the Cart sheet ---------------------------------------------
$wsCart = $wb.Worksheets.Item("Cart and Order Data Retrieval")
$wsCart.Name = "Cart Data Retrieval"

# --- Update the "User Data Retrieval" test case titles ------------------
$wsUser = $wb.Worksheets.Item("User Data Retrieval")
$wsUser.Range("B3").Value = "Fetch Valid Single User"
$wsUser.Range("B4").Value = "Fetch Invalid Single User"

# --- Update saved selections on the sheets that were visited ------------
$wsCart.Range("D22").Select()

$wsRecipe = $wb.Worksheets.Item("Recipe Data Retrieval")
$wsRecipe.Range("D31").Select()

# Leave "User Data Retrieval" as the active sheet with J9 selected, which
# becomes the workbook's activeTab on save.
$wsUser.Activate()
$wsUser.Range("J9").Select()
